$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain text so numeric-looking values
# (e.g. "25.552.91", "0.9979") are not reinterpreted as numbers/dates.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "25.552.91"
$ws.Range("E2").Value = "  +2.27%  "
$ws.Range("D3").Value = "1.664.75"
$ws.Range("E3").Value = "  +1.40%  "
$ws.Range("D4").Value = "0.9979"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "236.30"
$ws.Range("E5").Value = "  +1.52%  "
$ws.Range("D6").Value = "0.9988"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").Value = "0.4641"
$ws.Range("E7").Value = "  -2.51%  "
$ws.Range("D8").Value = "0.2574"
$ws.Range("E8").Value = "  -0.89%  "
$ws.Range("D9").Value = "0.06134"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").Value = "1.660.45"
$ws.Range("E10").Value = "  +1.13%  "
$ws.Range("D11").Value = "0.06942"
$ws.Range("E11").Value = "  -1.27%  "
$ws.Range("D12").Value = "14.80"
$ws.Range("E12").Value = "  +1.52%  "
$ws.Range("D13").Value = "4.331"
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").Value = "75.07"
$ws.Range("E14").Value = "  +1.94%  "
$ws.Range("D15").Value = "0.5712"
$ws.Range("E15").Value = "  -3.30%  "
$ws.Range("D16").Value = "0.9990"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").Value = "0.9993"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").Value = "25.539.23"
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("D19").Value = "0.000006701"
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("D20").Value = "11.38"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").Value = "1.873.12"
$ws.Range("E21").Value = "  +0.84%  "
$ws.Range("D22").Value = "4.420"
$ws.Range("E22").Value = "  +3.07%  "
$ws.Range("D23").Value = "8.612"
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("D24").Value = "5.225"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").Value = "134.00"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").Value = "14.93"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").Value = "1.371"
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("D28").Value = "1.714"
$ws.Range("E28").Value = "  +4.81%  "
$ws.Range("D29").Value = "104.08"
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("D30").Value = "3.939"
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("D31").Value = "0.07688"
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("D32").Value = "3.595"
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("D33").Value = "0.04353"
$ws.Range("E33").Value = "  +1.61%  "
$ws.Range("D34").Value = "2.604"
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("D35").Value = "0.6045"
$ws.Range("E35").Value = "  +2.37%  "
$ws.Range("D36").Value = "0.9382"
$ws.Range("E36").Value = "  +1.35%  "
$ws.Range("D37").Value = "0.9151"
$ws.Range("E37").Value = "  +3.28%  "
$ws.Range("D38").Value = "107.78"
$ws.Range("E38").Value = "  +8.87%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.366"
$ws.Range("E39").Value = "  -8.17%  "
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").Value = "0.9978"
$ws.Range("E40").Value = "  -0.23%  "
$ws.Range("D41").Value = "1.831"
$ws.Range("E41").Value = "  +4.43%  "
$ws.Range("D42").Value = "0.01449"
$ws.Range("E42").Value = "  -3.98%  "
$ws.Range("D43").Value = "0.3705"
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("D44").Value = "4.997"
$ws.Range("E44").Value = "  +7.01%  "
$ws.Range("D45").Value = "0.1108"
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("D46").Value = "0.05263"
$ws.Range("E46").Value = "  +1.07%  "
$ws.Range("D47").Value = "6.108"
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("D48").Value = "30.50"
$ws.Range("E48").Value = "  +5.55%  "
$ws.Range("D49").Value = "7.597"
$ws.Range("E49").Value = "  +6.91%  "
$ws.Range("D50").Value = "1.001"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").Value = "0.9979"
$ws.Range("E51").Value = "  -0.05%  "
